$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 2.9
$ws.Range("H2").Value = 2.75
$ws.Range("J2").Value = 4
$ws.Range("K2").Value = 1.73
$ws.Range("N2").Value = 4.75
$ws.Range("W2").Value = 5.5
$ws.Range("X2").Value = 12
$ws.Range("Z2").Value = 34
$ws.Range("AG2").Value = 5.5
$ws.Range("AO2").Value = 21
$ws.Range("AR2").Value = 151
$ws.Range("M4").Value = 1.07
$ws.Range("N4").Value = 9
$ws.Range("Q8").Value = 2.2
$ws.Range("R8").Value = 1.65
$ws.Range("Q10").Value = 1.62
$ws.Range("R10").Value = 2.25
$ws.Range("AD10").Value = 7.5
$ws.Range("AH10").Value = 23
$ws.Range("BA10").Value = 81
$ws.Range("G11").Value = 3.4
$ws.Range("H11").Value = 3
$ws.Range("I11").Value = 2.35
$ws.Range("L11").Value = 3.2
$ws.Range("O11").Value = 1.5
$ws.Range("P11").Value = 2.5
$ws.Range("Q11").Value = 2.6
$ws.Range("R11").Value = 1.48
$ws.Range("S11").Value = 1.57
$ws.Range("T11").Value = 2.25
$ws.Range("X11").Value = 15
$ws.Range("Z11").Value = 41
$ws.Range("AF11").Value = 67
$ws.Range("AH11").Value = 9.5
$ws.Range("AI11").Value = 10
$ws.Range("AJ11").Value = 21
$ws.Range("AO11").Value = 21
$ws.Range("AR11").Value = 126
$ws.Range("AT11").Value = 2.25
$ws.Range("M12").Value = 1.1
$ws.Range("N12").Value = 7
$ws.Range("Q12").Value = 2.5
$ws.Range("R12").Value = 1.5
$ws.Range("O13").Value = 1.67
$ws.Range("P13").Value = 2.1
